$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219, shifting existing rows 219-315 down to 220-316
$ws.Rows("219:219").Insert()

# Populate the new row 219 with the new data point
$ws.Cells.Item(219, 1).Value = 5
$ws.Cells.Item(219, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(219, 3).Value = "Maule"
$ws.Cells.Item(219, 4).Value = 44510
$ws.Cells.Item(219, 5).Value = 7
$ws.Cells.Item(219, 6).Value = 100114001
$ws.Cells.Item(219, 7).Value = "Papa"
$ws.Cells.Item(219, 8).Value = "Rodeo"
$ws.Cells.Item(219, 9).Value = "1a nueva(o)"
$ws.Cells.Item(219, 10).Value = 1600
$ws.Cells.Item(219, 11).Value = 9000
$ws.Cells.Item(219, 12).Value = 9000
$ws.Cells.Item(219, 13).Value = 9000
$ws.Cells.Item(219, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(219, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(219, 16).Value = 360
$ws.Cells.Item(219, 17).Value = 25
$ws.Cells.Item(219, 18).Value = "Hortaliza"
